$wb = $excel.ActiveWorkbook

# --- Sheet "LoginPageTest" ---
$wsLogin = $wb.Worksheets.Item("LoginPageTest")
$wsLogin.Select()
$wsLogin.Range("D7").Select()

# --- Sheet "TestSuite" ---
$wsSuite = $wb.Worksheets.Item("TestSuite")
$wsSuite.Range("B2").Value = "Y"
$wsSuite.Range("B3").Value = "Y"
$wsSuite.Range("B4").Value = "N"
$wsSuite.Select()
$wsSuite.Range("B2").Select()
